# Append a fresh scrape run (2025-12-08 06:31 JST) to the "ランサーズ" sheet.
# The whole data block (rows 2..N) is the latest top-N ranking snapshot, so we
# overwrite rows 2-8 with the new snapshot's first 7 entries and append two
# brand new rows (9-10) for the rest of the snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Capture the existing "Hyperlink" cell style (column F already uses it)
# before we touch anything, then drop the old hyperlinks outright - every
# row's URL is being rewritten (and two rows are brand new), so it's
# simplest to rebuild the hyperlink list from scratch rather than risk
# stacking duplicates on top of stale ones.
$linkStyle = $ws.Cells.Item(2, 6).Style
$ws.Hyperlinks.Delete()

$timestamp = "2025-12-08 06:31:34"

$rows = @(
    @{ Row = 2;  B = "【完全在宅】ChatGPT・AI活用講師募集|IT/業務支援経験者歓迎!"; D = "10,000 円 ~ 20,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5449394"; G = 600; H = "🔥AI,GPT" },
    @{ Row = 3;  B = "DreamWeaver – 夢日記 + 睡眠導入 + AI分析のアプリ開発"; D = "1,000 ~ 5,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5449048"; G = 370; H = "🔥AI,Ai ◆開発 ◇アプリ" },
    @{ Row = 4;  B = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"; D = "300,000 円 ~ 500,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5427956"; G = 310; H = "🔥AI,Ai" },
    @{ Row = 5;  B = "【フリーランス募集】訪問看護向けスケジュール管理アプリ開発"; D = "1,000,000 円 ~ 3,000,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5449280"; G = 135; H = "◆開発 ◇アプリ" },
    @{ Row = 6;  B = "【急募】紙の伝票をWEBシステムへ自動データ入力開発"; D = "300,000 円 ~ 500,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5449142"; G = 90;  H = "◆開発" },
    @{ Row = 7;  B = "自動出品システムの開発"; D = "100,000 円 ~ 200,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5449232"; G = 83;  H = "◆開発" },
    @{ Row = 8;  B = "初回 WebアプリのiOSアプリ化+IAPサブスク(2週無料)+申請"; D = "500,000 円 ~ 1,000,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5449067"; G = 45;  H = "◇アプリ" },
    @{ Row = 9;  B = "【急募】Shopifyでの3Dカスタムシミュレーター導入設定依頼"; D = "200,000 円 ~ 300,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5449335"; G = 18;  H = $null },
    @{ Row = 10; B = "初回 【継続案件】AWS上でのLAMP環境構築および保守・運用サポートパートナー募集"; D = "20,000 円 ~ 50,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5449313"; G = 13;  H = $null }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $timestamp
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = "システム開発"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = "期限情報なし"

    $urlCell = $ws.Cells.Item($row, 6)
    $urlCell.Value = $r.Url
    $ws.Hyperlinks.Add($urlCell, $r.Url)
    $urlCell.Style = $linkStyle

    $ws.Cells.Item($row, 7).Value = $r.G

    if ($r.H) {
        $ws.Cells.Item($row, 8).Value = $r.H
    }
}

$ws.Columns.Item(4).ColumnWidth = 31.166666666666668
